# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a completed
# handback: the Overview status text changes, the zh-cn report gets its
# handback datetime filled in, and the de-de report gets a full set of
# handback info (target file link, handback file, handback datetime).

$wb = $excel.ActiveWorkbook

$targetMdFile = "188b5bd1-3611-4d30-b56e-41c5a83314fc.md"
$targetMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25ec80ad937faebebf9cc690616d2256301e3021/e2e/188b5bd1-3611-4d30-b56e-41c5a83314fc.md"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Status text for both locales moves from "Ready for handoff" to
# "Handed back: in sync with en-US".
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# The status column got noticeably wider to fit the new text.
$wsOverview.Range("E1").ColumnWidth = 29.165
$wsOverview.Range("F1").ColumnWidth = 29.165

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C1").ColumnWidth = 29.165

# Latest Target File: link to the source md file, styled like the
# existing hyperlink in column A.
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $targetMdUrl, "", "", $targetMdFile)
$wsZhCn.Range("I1").ColumnWidth = 39.165

# Latest Handback File / DateTime
$wsZhCn.Range("J2").Value = "188b5bd1-3611-4d30-b56e-41c5a83314fc.302694c4ac57687c06ed79b1546a2c0aa0860d0a.zh-cn.xlf"
$wsZhCn.Range("J1").ColumnWidth = 39.165
$wsZhCn.Range("K2").Value = "2016-08-21 09:03:43"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C1").ColumnWidth = 29.165

# Latest Target File: link to the source md file, styled like the
# existing hyperlink in column A.
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $targetMdUrl, "", "", $targetMdFile)
$wsDeDe.Range("I1").ColumnWidth = 39.165

# Latest Handback File / DateTime
$wsDeDe.Range("J2").Value = "188b5bd1-3611-4d30-b56e-41c5a83314fc.302694c4ac57687c06ed79b1546a2c0aa0860d0a.de-de.xlf"
$wsDeDe.Range("J1").ColumnWidth = 39.165
$wsDeDe.Range("K2").Value = "2016-08-21 09:03:49"

Write-Host "Handback report generated"
